$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 271, shifting existing rows 271-292 down to 273-294
$ws.Rows.Item(271).Resize(2).Insert()

# New row 271 data
$ws.Cells.Item(271, 1).Value = 3
$ws.Cells.Item(271, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(271, 3).Value = "Coquimbo"
$ws.Cells.Item(271, 4).Value = 44578
$ws.Cells.Item(271, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(271, 5).Value = 5
$ws.Cells.Item(271, 6).Value = 100112031
$ws.Cells.Item(271, 7).Value = "Poroto verde"
$ws.Cells.Item(271, 8).Value = "Magnum"
$ws.Cells.Item(271, 9).Value = "Primera"
$ws.Cells.Item(271, 10).Value = 73
$ws.Cells.Item(271, 11).Value = 27000
$ws.Cells.Item(271, 12).Value = 28000
$ws.Cells.Item(271, 13).Value = 27479
$ws.Cells.Item(271, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(271, 15).Value = "Provincia de Talca"
$ws.Cells.Item(271, 16).Value = 1099
$ws.Cells.Item(271, 17).Value = 25
$ws.Cells.Item(271, 18).Value = "Hortaliza"

# New row 272 data
$ws.Cells.Item(272, 1).Value = 3
$ws.Cells.Item(272, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(272, 3).Value = "Coquimbo"
$ws.Cells.Item(272, 4).Value = 44578
$ws.Cells.Item(272, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(272, 5).Value = 5
$ws.Cells.Item(272, 6).Value = 100112031
$ws.Cells.Item(272, 7).Value = "Poroto verde"
$ws.Cells.Item(272, 8).Value = "Magnum"
$ws.Cells.Item(272, 9).Value = "Segunda"
$ws.Cells.Item(272, 10).Value = 38
$ws.Cells.Item(272, 11).Value = 21000
$ws.Cells.Item(272, 12).Value = 21000
$ws.Cells.Item(272, 13).Value = 21000
$ws.Cells.Item(272, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(272, 15).Value = "Provincia de Talca"
$ws.Cells.Item(272, 16).Value = 840
$ws.Cells.Item(272, 17).Value = 25
$ws.Cells.Item(272, 18).Value = "Hortaliza"
